# Updates cryptos list price (D) and volume (E) columns per the source diff.
# D2:D51 is temporarily forced to Text format before writing so that
# numeric-looking price strings (e.g. "246.42") are stored as text,
# matching the original inlineStr cells; the format/style is reset to
# Normal afterwards so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("D2").Value = '42.368.79'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '2.244.87'
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '246.42'
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("D6").Value = '0.631'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").Value = '76.05'
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '0.623'
$ws.Range("E9").Value = '  -1.56%  '
$ws.Range("D10").Value = '44.03'
$ws.Range("E10").Value = '  +9.88%  '
$ws.Range("D11").Value = '0.0950'
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").Value = '7.30'
$ws.Range("E12").Value = '  +0.92%  '
$ws.Range("E13").Value = '  -1.15%  '
$ws.Range("D14").Value = '2.588.75'
$ws.Range("E14").Value = '  +0.17%  '
$ws.Range("D15").Value = '14.61'
$ws.Range("E15").Value = '  -1.86%  '
$ws.Range("D16").Value = '0.857'
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("D17").Value = '2.248.59'
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("D18").Value = '42.203.53'
$ws.Range("E18").Value = '  -0.04%  '
$ws.Range("E19").Value = '  +4.11%  '
$ws.Range("D20").Value = '6.19'
$ws.Range("E20").Value = '  +0.65%  '
$ws.Range("D21").Value = '72.25'
$ws.Range("E22").Value = '  +2.73%  '
$ws.Range("D23").Value = '231.57'
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  +33.02%  '
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = '11.43'
$ws.Range("E26").Value = '  +2.48%  '
$ws.Range("E27").Value = '  -2.86%  '
$ws.Range("D28").Value = '2.32'
$ws.Range("E28").Value = '  -0.16%  '
$ws.Range("D29").Value = '2.19'
$ws.Range("E29").Value = '  +1.24%  '
$ws.Range("D30").Value = '168.09'
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("D31").Value = '20.71'
$ws.Range("E31").Value = '  +0.82%  '
$ws.Range("D32").Value = '0.0830'
$ws.Range("E32").Value = '  -2.69%  '
$ws.Range("E33").Value = '  +0.69%  '
$ws.Range("D34").Value = '30.70'
$ws.Range("E34").Value = '  -0.80%  '
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("E36").Value = '  +11.38%  '
$ws.Range("D37").Value = '4.54'
$ws.Range("E37").Value = '  +1.38%  '
$ws.Range("E38").Value = '  +6.94%  '
$ws.Range("D39").Value = '13.63'
$ws.Range("E39").Value = '  +5.04%  '
$ws.Range("D40").Value = '2.19'
$ws.Range("E40").Value = '  -1.90%  '
$ws.Range("E41").Value = '  -1.76%  '
$ws.Range("D42").Value = '63.54'
$ws.Range("E42").Value = '  +5.70%  '
$ws.Range("D43").Value = '0.203'
$ws.Range("E43").Value = '  -0.54%  '
$ws.Range("D44").Value = '108.16'
$ws.Range("E44").Value = '  -8.35%  '
$ws.Range("D45").Value = '8.79'
$ws.Range("E45").Value = '  +0.52%  '
$ws.Range("E46").Value = '  +1.51%  '
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("E48").Value = '  +1.53%  '
$ws.Range("E49").Value = '  -0.13%  '
$ws.Range("D50").Value = '2.35'
$ws.Range("E50").Value = '  +6.19%  '
$ws.Range("D51").Value = '4.11'
$ws.Range("E51").Value = '  -1.30%  '
$ws.Range("D2:D51").Style = "Normal"
Write-Output "Updated cryptos price/volume data."
